# Fixed update to excel issue
# 1) Rename header labels on the existing sheets to the new PO_Qty naming
$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Range("B1").Value = "Weekly_PO_Qty"

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add a new "PO Forecast" sheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$forecast = $wb.Worksheets.Add($null, $lastSheet)
$forecast.Name = "PO Forecast"

# Match the sheetPr outline settings used by the other sheets
$forecast.Outline.SummaryRow = 1
$forecast.Outline.SummaryColumn = 1

# Match the page margins used by the other sheets (0.75in/0.75in/1in/1in/0.5in/0.5in)
$forecast.PageSetup.LeftMargin = 54
$forecast.PageSetup.RightMargin = 54
$forecast.PageSetup.TopMargin = 72
$forecast.PageSetup.BottomMargin = 72
$forecast.PageSetup.HeaderMargin = 36
$forecast.PageSetup.FooterMargin = 36

# Header row, styled like the header rows on the other sheets
$weekly.Range("A1:B1").Copy()
$forecast.Range("A1:D1").PasteSpecial(-4122)

$forecast.Range("A1").Value = "ds"
$forecast.Range("B1").Value = "PO_Forecast"
$forecast.Range("C1").Value = "yhat_lower"
$forecast.Range("D1").Value = "yhat_upper"

# Date-formatted column A, matching the other sheets' date style
$weekly.Range("A2").Copy()
$forecast.Range("A2:A31").PasteSpecial(-4122)

$poData = @(
  @(45144.99999999999, 0, -309.9584288973914, 68.92008458619448),
  @(45158.99999999999, 0, -275.9740200112391, 97.66240340524286),
  @(45179.99999999999, 0, -229.0277570744432, 165.2299862578085),
  @(45186.99999999999, 0, -211.3528519150798, 190.6244962164252),
  @(45193.99999999999, 2, -173.5818711903755, 196.2225772599766),
  @(45207.99999999999, 37, -142.2583021282564, 233.8213240571929),
  @(45214.99999999999, 55, -145.1462777859865, 248.7854722588258),
  @(45221.99999999999, 72, -117.783709228381, 272.3157225653055),
  @(45228.99999999999, 89, -118.8520710332161, 288.2954293695909),
  @(45235.99999999999, 107, -84.10174747118805, 303.6148072029151),
  @(45242.99999999999, 124, -69.08082131303208, 332.0223789788963),
  @(45249.99999999999, 141, -52.06535693598829, 328.2133336145876),
  @(45256.99999999999, 159, -40.56737382535163, 344.5844499477979),
  @(45263.99999999999, 176, -16.17638286815213, 381.5578658783736),
  @(45270.99999999999, 194, -0.9693886141866439, 382.0537405440693),
  @(45277.99999999999, 211, 15.2334371512957, 401.9430564312699),
  @(45298.99999999999, 263, 61.94810443160678, 441.9491433847377),
  @(45305.99999999999, 280, 85.36674805736634, 467.6508674003221),
  @(45312.99999999999, 298, 106.040970514938, 497.3524017239074),
  @(45319.99999999999, 315, 125.2958251400976, 497.1012146686867),
  @(45326.99999999999, 333, 136.8942878736851, 529.8019970701837),
  @(45340.99999999999, 367, 159.7653382007402, 538.8035138180769),
  @(45347.99999999999, 385, 194.148921151251, 570.1382096110258),
  @(45354.99999999999, 402, 220.7895404959177, 593.0869404283529),
  @(45361.99999999999, 419, 227.9516268821707, 621.652109240131),
  @(45368.99999999999, 437, 246.162220618463, 634.0897970182838),
  @(45375.99999999999, 454, 262.9282132254041, 658.9858529359283),
  @(45382.99999999999, 472, 284.3420436174326, 666.820063705129),
  @(45389.99999999999, 489, 301.9093140705801, 678.8686113339016),
  @(45396.99999999999, 506, 322.1181501652826, 694.3030313005524)
)

$r = 2
foreach ($row in $poData) {
  $forecast.Cells.Item($r, 1).Value = $row[0]
  $forecast.Cells.Item($r, 2).Value = $row[1]
  $forecast.Cells.Item($r, 3).Value = $row[2]
  $forecast.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
